# Update countries & provincias Spain
# Applies a daily data refresh to the "Pais" sheet:
#  - Updates statistics for several existing countries (rows unchanged).
#  - Re-positions "Birmania" to immediately follow "Sierra Leona" (new stats),
#    shifting "Liberia", "Cabo Verde" and "Guadalupe" down by one row each
#    (their own data is carried over unchanged).
#  - Swaps the rows of "Burundi" and "San Cristobal y Nieves" so that
#    "Burundi" now comes immediately after "Curazao".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------
# 1) Straightforward statistic refreshes (country stays on its row)
# ---------------------------------------------------------------

# Estados Unidos (row 4)
$ws.Range("B4").Value = 1162383
$ws.Range("C4").Value = 1609
$ws.Range("E4").Value = 920968
$ws.Range("G4").Value = 61
$ws.Range("H4").Value = 67505

# Argentina (row 57)
$ws.Range("D57").Value = 1354
$ws.Range("E57").Value = 3086
$ws.Range("G57").Value = 4
$ws.Range("H57").Value = 241

# Sri Lanka (row 103)
$ws.Range("B103").Value = 706
$ws.Range("C103").Value = 4
$ws.Range("E103").Value = 517

# Maldivas (row 114)
$ws.Range("B114").Value = 527
$ws.Range("C114").Value = 8
$ws.Range("E114").Value = 508

# ---------------------------------------------------------------
# 2) Birmania moves up to right after Sierra Leona (row 138),
#    pushing Liberia / Cabo Verde / Guadalupe down one row each.
#    Capture the (unchanged) data of those three rows first.
# ---------------------------------------------------------------

$liberiaA = $ws.Range("A138").Value2
$liberiaB = $ws.Range("B138").Value2
$liberiaC = $ws.Range("C138").Value2
$liberiaD = $ws.Range("D138").Value2
$liberiaE = $ws.Range("E138").Value2
$liberiaF = $ws.Range("F138").Value2
$liberiaG = $ws.Range("G138").Value2
$liberiaH = $ws.Range("H138").Value2

$caboverdeA = $ws.Range("A139").Value2
$caboverdeB = $ws.Range("B139").Value2
$caboverdeC = $ws.Range("C139").Value2
$caboverdeD = $ws.Range("D139").Value2
$caboverdeE = $ws.Range("E139").Value2
$caboverdeF = $ws.Range("F139").Value2
$caboverdeG = $ws.Range("G139").Value2
$caboverdeH = $ws.Range("H139").Value2

$guadalupeA = $ws.Range("A140").Value2
$guadalupeB = $ws.Range("B140").Value2
$guadalupeC = $ws.Range("C140").Value2
$guadalupeD = $ws.Range("D140").Value2
$guadalupeE = $ws.Range("E140").Value2
$guadalupeF = $ws.Range("F140").Value2
$guadalupeG = $ws.Range("G140").Value2
$guadalupeH = $ws.Range("H140").Value2

# Row 138: Birmania, with its refreshed statistics
$ws.Range("A138").Value = "Birmania"
$ws.Range("B138").Value = 155
$ws.Range("C138").Value = 4
$ws.Range("D138").Value = 43
$ws.Range("E138").Value = 106
$ws.Range("F138").Value = 0
$ws.Range("G138").Value = 0
$ws.Range("H138").Value = 6

# Row 139: Liberia (carried over unchanged)
$ws.Range("A139").Value = $liberiaA
$ws.Range("B139").Value = $liberiaB
$ws.Range("C139").Value = $liberiaC
$ws.Range("D139").Value = $liberiaD
$ws.Range("E139").Value = $liberiaE
$ws.Range("F139").Value = $liberiaF
$ws.Range("G139").Value = $liberiaG
$ws.Range("H139").Value = $liberiaH

# Row 140: Cabo Verde (carried over unchanged)
$ws.Range("A140").Value = $caboverdeA
$ws.Range("B140").Value = $caboverdeB
$ws.Range("C140").Value = $caboverdeC
$ws.Range("D140").Value = $caboverdeD
$ws.Range("E140").Value = $caboverdeE
$ws.Range("F140").Value = $caboverdeF
$ws.Range("G140").Value = $caboverdeG
$ws.Range("H140").Value = $caboverdeH

# Row 141: Guadalupe (carried over unchanged)
$ws.Range("A141").Value = $guadalupeA
$ws.Range("B141").Value = $guadalupeB
$ws.Range("C141").Value = $guadalupeC
$ws.Range("D141").Value = $guadalupeD
$ws.Range("E141").Value = $guadalupeE
$ws.Range("F141").Value = $guadalupeF
$ws.Range("G141").Value = $guadalupeG
$ws.Range("H141").Value = $guadalupeH

# ---------------------------------------------------------------
# 3) Burundi moves up to right after Curazao (row 198), swapping
#    places with San Cristobal y Nieves (row 199). Both rows keep
#    their own data unchanged, only their positions are swapped.
# ---------------------------------------------------------------

$sancrisA = $ws.Range("A198").Value2
$sancrisB = $ws.Range("B198").Value2
$sancrisC = $ws.Range("C198").Value2
$sancrisD = $ws.Range("D198").Value2
$sancrisE = $ws.Range("E198").Value2
$sancrisF = $ws.Range("F198").Value2
$sancrisG = $ws.Range("G198").Value2
$sancrisH = $ws.Range("H198").Value2

$burundiA = $ws.Range("A199").Value2
$burundiB = $ws.Range("B199").Value2
$burundiC = $ws.Range("C199").Value2
$burundiD = $ws.Range("D199").Value2
$burundiE = $ws.Range("E199").Value2
$burundiF = $ws.Range("F199").Value2
$burundiG = $ws.Range("G199").Value2
$burundiH = $ws.Range("H199").Value2

# Row 198: Burundi
$ws.Range("A198").Value = $burundiA
$ws.Range("B198").Value = $burundiB
$ws.Range("C198").Value = $burundiC
$ws.Range("D198").Value = $burundiD
$ws.Range("E198").Value = $burundiE
$ws.Range("F198").Value = $burundiF
$ws.Range("G198").Value = $burundiG
$ws.Range("H198").Value = $burundiH

# Row 199: San Cristobal y Nieves
$ws.Range("A199").Value = $sancrisA
$ws.Range("B199").Value = $sancrisB
$ws.Range("C199").Value = $sancrisC
$ws.Range("D199").Value = $sancrisD
$ws.Range("E199").Value = $sancrisE
$ws.Range("F199").Value = $sancrisF
$ws.Range("G199").Value = $sancrisG
$ws.Range("H199").Value = $sancrisH
